# Change the header title from "Assembly Guide" to "mAKER Guide"
# (displayed as "MAKER GUIDE" because of the w:caps run formatting).
# The original text is a single run; the edit splits it into two runs:
#   "mAKER" and " Guide" (leading space preserved), both carrying the
# same run formatting (Roboto, bold, all-caps, color 646464, size 32).

$d = $word.ActiveDocument

# Locate the header paragraph containing "Assembly Guide".
$hdr = $null
foreach ($sec in $d.Sections) {
    foreach ($candidate in $sec.Headers) {
        $probe = $candidate.Range.Duplicate
        if ($probe.Find.Execute("Assembly Guide")) {
            $hdr = $candidate
            break
        }
    }
    if ($hdr -ne $null) { break }
}

$rng = $hdr.Range.Duplicate
$found = $rng.Find.Execute("Assembly Guide")

# Clear the existing run's text, then inject the two replacement runs
# (with explicit run properties matching the original formatting) via
# InsertXML so the paragraph mark / paragraph properties are untouched
# and the run split + per-run formatting is reproduced exactly.
$rng.Text = ""

$xmlPkg = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="37BB2FA2" w14:textId="1B255A56" w:rsidR="16EBA8AB" w:rsidRDefault="4CE3AF1C" w:rsidP="4CE3AF1C"><w:pPr><w:pStyle w:val="Header"/><w:rPr><w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/><w:b/><w:bCs/><w:caps/><w:color w:val="646464"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:r w:rsidRPr="4CE3AF1C"><w:rPr><w:rFonts w:ascii="Roboto" w:eastAsia="Roboto" w:hAnsi="Roboto" w:cs="Roboto"/><w:b/><w:bCs/><w:caps/><w:color w:val="646464"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>mAKER</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Roboto" w:eastAsia="Roboto" w:hAnsi="Roboto" w:cs="Roboto"/><w:b/><w:bCs/><w:caps/><w:color w:val="646464"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve"> Guide</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$rng.InsertXML($xmlPkg)

Write-Output $hdr.Range.Text
